$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'52.175.74"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "'2.935.66"
$ws.Range("E3").Value = "  +4.29%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'352.94"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "'113.40"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "'0.559"
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.623"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").Value = "'39.61"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("D11").Value = "'0.0880"
$ws.Range("E11").Value = "  +4.09%  "
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").Value = "'20.06"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "'7.77"
$ws.Range("E14").Value = "  -0.47%  "
$ws.Range("D15").Value = "'3.389.62"
$ws.Range("E15").Value = "  +3.93%  "
$ws.Range("D16").Value = "'2.926.11"
$ws.Range("E16").Value = "  +4.11%  "
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Value = "'52.191.62"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "'7.62"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").Value = "'3.30"
$ws.Range("E20").Value = "  -2.38%  "
$ws.Range("D21").Value = "'14.20"
$ws.Range("E21").Value = "  +4.37%  "
$ws.Range("D22").Value = "'0.0₃0980"
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("D23").Value = "'71.16"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").Value = "'269.19"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  +0.66%  "
$ws.Range("D26").Value = "'0.179"
$ws.Range("E26").Value = "  +9.57%  "
$ws.Range("D27").Value = "'27.05"
$ws.Range("E27").Value = "  +2.98%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").Value = "'6.99"
$ws.Range("E29").Value = "  +13.23%  "
$ws.Range("D30").Value = "'10.64"
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").Value = "'0.103"
$ws.Range("E31").Value = "  +14.11%  "
$ws.Range("D32").Value = "'2.27"
$ws.Range("E32").Value = "  -0.42%  "
$ws.Range("D33").Value = "'37.15"
$ws.Range("E33").Value = "  -4.42%  "
$ws.Range("D34").Value = "'6.05"
$ws.Range("E34").Value = "  +5.75%  "
$ws.Range("D35").Value = "'53.09"
$ws.Range("E35").Value = "  +0.42%  "
$ws.Range("D36").Value = "'0.0454"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("E37").Value = "  -0.24%  "
$ws.Range("D38").Value = "'3.38"
$ws.Range("E38").Value = "  +5.24%  "
$ws.Range("D39").Value = "'18.69"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("D40").Value = "'2.05"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").Value = "'2.70"
$ws.Range("E41").Value = "  +3.93%  "
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("D43").Value = "'23.06"
$ws.Range("E43").Value = "  +2.83%  "
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("D45").Value = "'2.192.58"
$ws.Range("E45").Value = "  +2.55%  "
$ws.Range("D46").Value = "'3.52"
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").Value = "'111.26"
$ws.Range("E48").Value = "  -8.23%  "
$ws.Range("D49").Value = "'0.249"
$ws.Range("E49").Value = "  +10.72%  "
$ws.Range("E50").Value = "  +7.00%  "
$ws.Range("D51").Value = "'0.955"
$ws.Range("E51").Value = "  -7.10%  "
